$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price-column values now look like plain numbers
# (e.g. "1.013", "309.41"). Assigning such strings via .Value would make
# Excel auto-convert them to actual numbers, which would not match the
# original text-based storage used throughout this sheet. Force those
# specific cells to Text format first so the values stick as strings,
# exactly like the rest of the column.
$numericTextCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D21", "D22", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.969.03'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '1.845.41'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  +0.41%  '
$ws.Range("D5").Value = '1.013'
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").Value = '309.41'
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("E7").Value = '  +1.92%  '
$ws.Range("D8").Value = '0.3681'
$ws.Range("E8").Value = '  +1.23%  '
$ws.Range("D9").Value = '0.07234'
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").Value = '0.9317'
$ws.Range("E10").Value = '  -0.28%  '
$ws.Range("D11").Value = '19.83'
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").Value = '0.07744'
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").Value = '1.876.29'
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("D14").Value = '5.384'
$ws.Range("E14").Value = '  +1.61%  '
$ws.Range("D15").Value = '6.468'
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").Value = '88.90'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '1.016'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").Value = '0.000008668'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").Value = '27.023.02'
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").Value = '14.57'
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("D22").Value = '5.075'
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").Value = '  +0.85%  '
$ws.Range("D25").Value = '153.08'
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("D26").Value = '18.24'
$ws.Range("E26").Value = '  +1.06%  '
$ws.Range("D27").Value = '2.015'
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("D28").Value = '114.44'
$ws.Range("E28").Value = '  +0.32%  '
$ws.Range("D29").Value = '4.977'
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("D30").Value = '0.08873'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").Value = '3.311'
$ws.Range("E31").Value = '  +4.02%  '
$ws.Range("D32").Value = '1.181'
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").Value = '0.7431'
$ws.Range("E33").Value = '  -0.49%  '
$ws.Range("D34").Value = '4.502'
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("D35").Value = '2.697'
$ws.Range("E35").Value = '  -5.74%  '
$ws.Range("D36").Value = '1.116'
$ws.Range("E36").Value = '  +2.22%  '
$ws.Range("D37").Value = '0.01963'
$ws.Range("E37").Value = '  +1.22%  '
$ws.Range("D38").Value = '0.05265'
$ws.Range("E38").Value = '  +1.93%  '
$ws.Range("D39").Value = '2.970'
$ws.Range("E39").Value = '  -0.66%  '
$ws.Range("D40").Value = '0.5278'
$ws.Range("E40").Value = '  +2.74%  '
$ws.Range("D41").Value = '7.022'
$ws.Range("E41").Value = '  +1.64%  '
$ws.Range("D42").Value = '0.1512'
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("D43").Value = '8.295'
$ws.Range("E43").Value = '  +1.29%  '
$ws.Range("D44").Value = '10.55'
$ws.Range("E44").Value = '  +0.88%  '
$ws.Range("D45").Value = '0.4741'
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("D46").Value = '1.015'
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("D47").Value = '102.02'
$ws.Range("E47").Value = '  +1.70%  '
$ws.Range("D48").Value = '1.608'
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").Value = '66.02'
$ws.Range("D50").Value = '0.06077'
$ws.Range("E50").Value = '  +0.50%  '
$ws.Range("D51").Value = '0.8922'
$ws.Range("E51").Value = '  +3.56%  '
Write-Host "Updated 95 Price/Volume cells in cryptos sheet."
